$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.912.58"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.776.96"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.29"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5371"
$ws.Range("E7").Value = "  -2.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3758"
$ws.Range("E8").Value = "  -2.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07431"
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.66"
$ws.Range("E10").Value = "  -2.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.089"
$ws.Range("E11").Value = "  -2.85%  "
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.38"
$ws.Range("E13").Value = "  -3.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.063"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.194"
$ws.Range("E15").Value = "  -2.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.769.29"
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.76"
$ws.Range("E17").Value = "  -4.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001051"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06402"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.865"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.949.08"
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.10"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.079"
$ws.Range("E25").Value = "  -2.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.77"
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.18"
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.975.69"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.268"
$ws.Range("E29").Value = "  -5.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.76"
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.109"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1046"
$ws.Range("E32").Value = "  +2.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.639"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.498"
$ws.Range("E34").Value = "  -4.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2232"
$ws.Range("E35").Value = "  -5.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06344"
$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02260"
$ws.Range("E37").Value = "  -2.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.947"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.374"
$ws.Range("E39").Value = "  -5.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6103"
$ws.Range("E40").Value = "  -4.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.177"
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.98"
$ws.Range("E42").Value = "  -5.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.427"
$ws.Range("E43").Value = "  +3.31%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.650"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5724"
$ws.Range("E47").Value = "  -4.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.66"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.180"
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.917"
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06781"
$ws.Range("E51").Value = "  -1.83%  "
